$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first sheet ("Hoja1" -> "Aulas")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Aulas"
$ws2 = $wb.Worksheets.Item("Profesores")

# ---------------------------------------------------------------------------
# 2. Wipe existing contents so stale shared-strings / formatting are dropped
# ---------------------------------------------------------------------------
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# ---------------------------------------------------------------------------
# 3. Re-populate "Profesores" (sheet2) first so the shared-string table gets
#    built in the same order as the target workbook.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = 1546
$ws2.Range("A1").NumberFormat = "@"

$ws2.Range("B1").Value = 1
$ws2.Range("B1").NumberFormat = "@"

$ws2.Range("C1").NumberFormat = "@"
$ws2.Range("C1").Value = "Fulanito Martinez"

$ws2.Range("D1").NumberFormat = "@"
$ws2.Range("D1").Value = "ZXC"

$ws2.Range("E1").NumberFormat = "@"
$ws2.Range("E1").Value = "BNM"

$ws2.Range("A2").Value = 6465
$ws2.Range("A2").NumberFormat = "@"

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "1"

$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "Arturo Perez Reverte"

$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "QWE"

$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "RTY"

# Column widths for Profesores
$ws2.Columns.Item(1).ColumnWidth = 4.166666666666667
$ws2.Columns.Item(2).ColumnWidth = 1.1666666666666665
$ws2.Columns.Item(3).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 4.276041666666667
$ws2.Columns.Item(5).ColumnWidth = 4.276041666666667

$ws2.Activate()
$ws2.Range("E2").Select()

# ---------------------------------------------------------------------------
# 4. Re-populate "Aulas" (sheet1)
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = 8484
$ws1.Range("A1").NumberFormat = "@"

$ws1.Range("B1").NumberFormat = "@"
$ws1.Range("B1").Value = "A201-B"

$ws1.Range("C1").NumberFormat = "@"
$ws1.Range("C1").Value = "Laboratorio"

$ws1.Range("D1").Value = 30
$ws1.Range("D1").NumberFormat = "@"

$ws1.Range("E1").NumberFormat = "@"
$ws1.Range("E1").Value = "idk"

$ws1.Range("F1").NumberFormat = "@"
$ws1.Range("F1").Value = "E201"

$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "7845"

$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "A202-C"

$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "Salon"

$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "30"

$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "ayuda"

$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("F2").Value = "me"

# ---------------------------------------------------------------------------
# 5. Add the new, empty "Hoja2" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Hoja2"

# ---------------------------------------------------------------------------
# 6. Leave "Aulas" as the active sheet/selection, matching the target view
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D12").Select()
